# Auto-generated Excel COM-interop script
# Applies refreshed market-price data (columns H-N) across all 8 class sheets
# as produced by the scheduled Faerie Profits price-update runner.
$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H9").Value = 1200715.6
$ws.Range("I9").Value = 1500647
$ws.Range("K9").Value = 1500647
$ws.Range("M9").Value = -1500478
$ws.Range("H17").Value = 3998.5173
$ws.Range("J17").Value = 4105.607
$ws.Range("L17").Value = 12316.821
$ws.Range("N17").Value = -12652.821
$ws.Range("H18").Value = 966.7143
$ws.Range("I18").Value = 966.7143
$ws.Range("K18").Value = 966.7143
$ws.Range("M18").Value = -682.7143
$ws.Range("H43").Value = 1431409.5
$ws.Range("J43").Value = 12767.125
$ws.Range("L43").Value = 12767.125
$ws.Range("N43").Value = -12905.125
$ws.Range("H112").Value = 3523.1636
$ws.Range("J112").Value = 3534.3462
$ws.Range("L112").Value = 10603.0386
$ws.Range("N112").Value = -12819.0386
$ws.Range("H129").Value = 3648.5217
$ws.Range("I129").Value = 851.3333
$ws.Range("K129").Value = 2553.9999
$ws.Range("M129").Value = 2446.0001
$ws.Range("H132").Value = 2091.077
$ws.Range("I132").Value = 2092.5
$ws.Range("K132").Value = 6277.5
$ws.Range("M132").Value = -3747.5
$ws.Range("H137").Value = 2922.4614
$ws.Range("I137").Value = 3330.2222
$ws.Range("J137").Value = 2005
$ws.Range("K137").Value = 9990.6666
$ws.Range("L137").Value = 6015
$ws.Range("M137").Value = -7440.6666
$ws.Range("N137").Value = -11115
$ws.Range("H138").Value = 142593.36
$ws.Range("I138").Value = 6061.778
$ws.Range("J138").Value = 161497.73
$ws.Range("K138").Value = 18185.334
$ws.Range("L138").Value = 484493.1900000001
$ws.Range("M138").Value = -13045.334
$ws.Range("N138").Value = -494773.1900000001

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 21810.488
$ws.Range("I32").Value = 18791.104
$ws.Range("K32").Value = 18791.104
$ws.Range("M32").Value = -18504.104
$ws.Range("H45").Value = 2637.923
$ws.Range("J45").Value = 3701.111
$ws.Range("L45").Value = 3701.111
$ws.Range("N45").Value = -4455.111
$ws.Range("H74").Value = 4084004.2
$ws.Range("I74").Value = 4763511.5
$ws.Range("K74").Value = 4763511.5
$ws.Range("M74").Value = -4762637.5
$ws.Range("H77").Value = 4084004.2
$ws.Range("I77").Value = 4763511.5
$ws.Range("K77").Value = 23817557.5
$ws.Range("M77").Value = -23813189.5
$ws.Range("H122").Value = 6184.448
$ws.Range("I122").Value = 5781.1665
$ws.Range("K122").Value = 17343.4995
$ws.Range("M122").Value = -14893.4995

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H134").Value = 3973.375
$ws.Range("I134").Value = 3739.6924
$ws.Range("K134").Value = 11219.0772
$ws.Range("M134").Value = -8684.0772

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H6").Value = 7600000
$ws.Range("I6").Value = 7600000
$ws.Range("K6").Value = 7600000
$ws.Range("M6").Value = -7599887
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 600
$ws.Range("K22").Value = 600
$ws.Range("M22").Value = -250
$ws.Range("H41").Value = 3000
$ws.Range("I41").Value = 3000
$ws.Range("K41").Value = 3000
$ws.Range("M41").Value = -2572
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H51").Value = 65000
$ws.Range("J51").Value = 65000
$ws.Range("L51").Value = 65000
$ws.Range("N51").Value = -66472
$ws.Range("H60").Value = 51666.332
$ws.Range("J60").Value = 57499.5
$ws.Range("L60").Value = 57499.5
$ws.Range("N60").Value = -58521.5
$ws.Range("H61").Value = 65000
$ws.Range("J61").Value = 65000
$ws.Range("L61").Value = 65000
$ws.Range("N61").Value = -65696
$ws.Range("H99").Value = 5999.4
$ws.Range("I99").Value = 5927.857
$ws.Range("K99").Value = 5927.857
$ws.Range("M99").Value = -4429.857
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("H126").Value = 5999.4
$ws.Range("I126").Value = 5927.857
$ws.Range("K126").Value = 17783.571
$ws.Range("M126").Value = -15313.571
$ws.Range("H132").Value = 7000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H134").Value = 3709.9443
$ws.Range("I134").Value = 1594.7222
$ws.Range("K134").Value = 4784.1666
$ws.Range("M134").Value = -2249.1666

# --- Sheet 5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H39").Value = 10982.071
$ws.Range("I39").Value = 937.25
$ws.Range("K39").Value = 2811.75
$ws.Range("M39").Value = -2517.75
$ws.Range("H56").Value = 7582.923
$ws.Range("I56").Value = 7582.923
$ws.Range("K56").Value = 7582.923
$ws.Range("M56").Value = -7052.923
$ws.Range("H107").Value = 3348.7307
$ws.Range("J107").Value = 3615.348
$ws.Range("L107").Value = 10846.044
$ws.Range("N107").Value = -14686.044
$ws.Range("H122").Value = 2479.6
$ws.Range("J122").Value = 3999
$ws.Range("L122").Value = 35991
$ws.Range("N122").Value = -40891
$ws.Range("H131").Value = 1431141.1
$ws.Range("J131").Value = 2998.1667
$ws.Range("L131").Value = 8994.500100000001
$ws.Range("N131").Value = -19074.5001

# --- Sheet 6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H57").Value = 21053.268
$ws.Range("J57").Value = 69993
$ws.Range("L57").Value = 69993
$ws.Range("N57").Value = -71633
$ws.Range("H126").Value = 45686.707
$ws.Range("I126").Value = 54004.9
$ws.Range("K126").Value = 162014.7
$ws.Range("M126").Value = -159544.7
$ws.Range("H132").Value = 4873.392
$ws.Range("I132").Value = 5051.136
$ws.Range("K132").Value = 15153.408
$ws.Range("M132").Value = -12623.408

# --- Sheet 7 ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 15725.286
$ws.Range("I7").Value = 17016
$ws.Range("J7").Value = 12498.5
$ws.Range("K7").Value = 17016
$ws.Range("L7").Value = 12498.5
$ws.Range("M7").Value = -16904
$ws.Range("N7").Value = -12722.5
$ws.Range("H22").Value = 2720.7856
$ws.Range("I22").Value = 2398.25
$ws.Range("K22").Value = 2398.25
$ws.Range("M22").Value = -2103.25
$ws.Range("H27").Value = 2720.7856
$ws.Range("I27").Value = 2398.25
$ws.Range("K27").Value = 2398.25
$ws.Range("M27").Value = -2291.25
$ws.Range("H46").Value = 5922.9395
$ws.Range("I46").Value = 3120
$ws.Range("K46").Value = 3120
$ws.Range("M46").Value = -2932
$ws.Range("H61").Value = 24056.88
$ws.Range("I61").Value = 32737.166
$ws.Range("K61").Value = 32737.166
$ws.Range("M61").Value = -32535.166
$ws.Range("H109").Value = 82500
$ws.Range("J109").Value = 82500
$ws.Range("L109").Value = 82500
$ws.Range("N109").Value = -85274
$ws.Range("H111").Value = 90000
$ws.Range("J111").Value = 90000
$ws.Range("L111").Value = 90000
$ws.Range("N111").Value = -98180
$ws.Range("H113").Value = 24056.88
$ws.Range("I113").Value = 32737.166
$ws.Range("K113").Value = 32737.166
$ws.Range("M113").Value = -30567.166
$ws.Range("H117").Value = 94995
$ws.Range("I117").Value = 94995
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 94995
$ws.Range("M117").Value = -90406
$ws.Range("N117").ClearContents()
$ws.Range("H126").Value = 15725.286
$ws.Range("I126").Value = 17016
$ws.Range("J126").Value = 12498.5
$ws.Range("K126").Value = 51048
$ws.Range("L126").Value = 37495.5
$ws.Range("M126").Value = -48578
$ws.Range("N126").Value = -42435.5
$ws.Range("H136").Value = 6499.8
$ws.Range("I136").Value = 5333.067
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 15999.201
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -13449.201
$ws.Range("N136").Value = -35100

# --- Sheet 8 ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H4").Value = 125013750
$ws.Range("I4").Value = 9333
$ws.Range("K4").Value = 9333
$ws.Range("M4").Value = -9220
$ws.Range("H99").Value = 45000
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H118").Value = 79995
$ws.Range("I118").Value = 79995
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 79995
$ws.Range("M118").Value = -78338
$ws.Range("N118").ClearContents()
$ws.Range("H126").Value = 2442.087
$ws.Range("I126").Value = 2448.0557
$ws.Range("J126").Value = 2420.6
$ws.Range("K126").Value = 7344.1671
$ws.Range("L126").Value = 7261.799999999999
$ws.Range("M126").Value = -4874.1671
$ws.Range("N126").Value = -12201.8
$ws.Range("H130").Value = 150000
$ws.Range("J130").Value = 150000
$ws.Range("L130").Value = 150000
$ws.Range("N130").Value = -160040
